$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# New data for rows 2-14 (one full cycle: TestScenario_1..4 covering
# New/View/Edit/Delete Account test cases). Rows 15-27 are an exact
# duplicate of this block (mirrors the pre-existing duplication
# pattern that was already present for TestScenario_1 in rows 2-12
# and 13-23 before the edit).
# -----------------------------------------------------------------
$data = @(
    @("TestScenario_1", "TestScenario_1.TestCase_1", "New Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab, and click on New button", "User should be navigated to the New  Account Page", "", ""),
    @("", "", "", "", "Valid value for required field Account Name ", "Step 2", "Input valid value in the  Account Name field.", "User should be able to input value for the Account Name field.", "", ""),
    @("", "", "", "", "Valid value for required field  ", "Step 3", "Input valid value in the   field.", "User should be able to input value for the  field.", "", ""),
    @("", "", "", "", "", "Step 4", "Click on Save button to save Account with fields", "User should be able to validate that a New Account is created", "", ""),
    @("TestScenario_2", "TestScenario_2.TestCase_1", "View Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab,  and select a Account ", "User should be navigated to the Account Page", "", ""),
    @("", "", "", "", "", "Step 2", "Click on Account name to View the Details", "User should be able to view the Account Details", "", ""),
    @("TestScenario_3", "TestScenario_3.TestCase_1", "Edit Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab,  and click on existing Account to modify", "User is navigated to the Account Details page", "", ""),
    @("", "", "", "", "Valid value for required field Account Name ", "Step 2", "Input valid value in the  Account Name field.", "User should be able to input value for the Account Name field.", "", ""),
    @("", "", "", "", "Valid value for required field  ", "Step 3", "Input valid value in the   field.", "User should be able to input value for the  field.", "", ""),
    @("", "", "", "", "", "Step 4", "Click on Save button to save Account with fields", "User should be able to validate that the Account is edited", "", ""),
    @("TestScenario_4", "TestScenario_4.TestCase_1", "Delete Account", "User Needs to Login to Salesforce, from the browser with correct credentials", "", "Step 1", "Click on the Account tab,  and select the existing  Account to delete", "User is navigated to the Account Details page", "", ""),
    @("", "", "", "", "", "Step 2", "Click on to the Delete to Delete the Account", "User should be able to validate that a pop-up is displayed asking for confirmation to delete the Account", "", ""),
    @("", "", "", "", "", "Step 3", "Click on Confirm / OK to delete the  Account", "User should be able to validate the Account is deleted", "", "")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  for ($j = 0; $j -lt 10; $j++) {
    $col = $j + 1
    $ws.Cells.Item($row, $col).Value = $data[$i][$j]
  }
}

# Duplicate the same 13-row block into rows 15-27
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 15
  for ($j = 0; $j -lt 10; $j++) {
    $col = $j + 1
    $ws.Cells.Item($row, $col).Value = $data[$i][$j]
  }
}

# -----------------------------------------------------------------
# Grow the table / autofilter range from A1:J23 to A1:J27
# -----------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J27"))

# -----------------------------------------------------------------
# Updated column widths for columns C-H
# -----------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 69.16666666666667
$ws.Columns.Item(5).ColumnWidth = 41.166666666666664
$ws.Columns.Item(6).ColumnWidth = 7.0
$ws.Columns.Item(7).ColumnWidth = 61.0
$ws.Columns.Item(8).ColumnWidth = 92.66666666666667
